$d = $word.ActiveDocument

# The edit simplifies the <w:docDefaults> block in word/styles.xml:
#   - rPrDefault/rPr keeps only rFonts, sz, szCs, lang (drops b, i,
#     smallCaps, strike, color, u, shd, vertAlign)
#   - pPrDefault/pPr keeps only a <w:spacing> element with just
#     line/lineRule (drops keepNext, keepLines, widowControl, pBdr, shd,
#     spacing's after/before, ind, contextualSpacing, jc)
#
# There's no Word object-model surface for docDefaults directly (setting
# Styles("Normal").Font/.ParagraphFormat creates direct overrides on the
# Normal style instead of touching docDefaults), so we round-trip the
# whole package through WordOpenXML and patch the docDefaults fragment
# with a literal string replace.

$xml = $d.WordOpenXML

$oldDocDefaults = '<w:docDefaults><w:rPrDefault><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b w:val="0"/><w:i w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/><w:shd w:val="clear" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:lang w:val="en"/></w:rPr></w:rPrDefault><w:pPrDefault><w:pPr><w:keepNext w:val="0"/><w:keepLines w:val="0"/><w:widowControl/><w:pBdr><w:top w:val="nil" w:sz="0" w:space="0"/><w:left w:val="nil" w:sz="0" w:space="0"/><w:bottom w:val="nil" w:sz="0" w:space="0"/><w:right w:val="nil" w:sz="0" w:space="0"/><w:between w:val="nil" w:sz="0" w:space="0"/></w:pBdr><w:shd w:val="clear" w:fill="auto"/><w:spacing w:before="0" w:after="0" w:line="276" w:lineRule="auto"/><w:ind w:left="0" w:right="0" w:firstLine="0"/><w:contextualSpacing w:val="0"/><w:jc w:val="left"/></w:pPr></w:pPrDefault></w:docDefaults>'

$newDocDefaults = '<w:docDefaults><w:rPrDefault><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en"/></w:rPr></w:rPrDefault><w:pPrDefault><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr></w:pPrDefault></w:docDefaults>'

if ($xml.IndexOf($oldDocDefaults) -lt 0) {
    throw "docDefaults fragment not found in WordOpenXML; cannot apply edit"
}

$xml = $xml.Replace($oldDocDefaults, $newDocDefaults)

$d.WordOpenXML = $xml
